# Updated cryptos list (price + 1h volume change refresh, plus an OKB/Toncoin
# row-order swap at rows 34-35) matching the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" on numeric-looking Price values forces Excel to keep them as
# text (matching the source inlineStr cells) instead of auto-converting to a
# number.
$ws.Range("D2").Value = "51.818.35"
$ws.Range("E2").Value = "  -0.98%  "

$ws.Range("D3").Value = "2.916.09"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'357.57"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("D6").Value = "'109.02"
$ws.Range("E6").Value = "  -3.52%  "

$ws.Range("D7").Value = "'0.564"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.629"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").Value = "'39.15"
$ws.Range("E10").Value = "  -3.06%  "

$ws.Range("D11").Value = "'0.0871"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("D13").Value = "'19.52"
$ws.Range("E13").Value = "  -3.05%  "

$ws.Range("D14").Value = "'7.82"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").Value = "3.383.00"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").Value = "2.926.11"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "'0.984"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").Value = "51.855.02"
$ws.Range("E18").Value = "  -1.02%  "

$ws.Range("D19").Value = "'3.34"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'7.55"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").Value = "'13.87"
$ws.Range("E21").Value = "  -4.35%  "

$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").Value = "'70.83"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").Value = "'269.46"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("D25").Value = "'2.82"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("E26").Value = "  +13.01%  "

$ws.Range("D27").Value = "'26.91"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").Value = "'7.56"
$ws.Range("E28").Value = "  +15.33%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  +12.17%  "

$ws.Range("D31").Value = "'10.54"
$ws.Range("E31").Value = "  -1.04%  "

$ws.Range("D32").Value = "'37.86"
$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("D33").Value = "'6.02"
$ws.Range("E33").Value = "  -2.42%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'52.40"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "'2.10"
$ws.Range("E35").Value = "  -7.01%  "

$ws.Range("D36").Value = "'0.0442"
$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").Value = "'3.22"
$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("D39").Value = "'18.26"
$ws.Range("E39").Value = "  -3.94%  "

$ws.Range("D40").Value = "'2.00"
$ws.Range("E40").Value = "  -4.14%  "

$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "  -4.07%  "

$ws.Range("D42").Value = "'0.119"
$ws.Range("E42").Value = "  +1.56%  "

$ws.Range("D43").Value = "'22.75"
$ws.Range("E43").Value = "  -5.09%  "

$ws.Range("D44").Value = "'119.48"
$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("D45").Value = "'2.18"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").Value = "'3.47"
$ws.Range("E46").Value = "  -2.72%  "

$ws.Range("E47").Value = "  -6.60%  "

$ws.Range("D48").Value = "2.122.72"
$ws.Range("E48").Value = "  -4.37%  "

$ws.Range("D49").Value = "'0.249"
$ws.Range("E49").Value = "  -4.76%  "

$ws.Range("D50").Value = "'0.0333"
$ws.Range("E50").Value = "  -2.01%  "

$ws.Range("D51").Value = "'9.13"
$ws.Range("E51").Value = "  -0.27%  "
